# "Small additions and final test" - update the disabilities roster so
# that the member previously listed as "Member E" is corrected/renamed
# to "Member D".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Disabilities")

# The member name lives in cell A5 of the Disabilities table.
$ws.Range("A5").Value = "Member D"

# Leave the freshly-edited row selected, as it was after the edit.
$ws.Range("A5").Select()
